$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header row: add "Wins", "Losses", "Ties" in AD1:AF1
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match formatting of the existing header cells (bold, centered, bordered)
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Data rows 2-40: team record values (Wins=103, Losses=58, Ties=0) for every player row
for ($r = 2; $r -le 40; $r++) {
    $ws.Cells.Item($r, 30).Value = 103
    $ws.Cells.Item($r, 31).Value = 58
    $ws.Cells.Item($r, 32).Value = 0
}
